$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0
$ws.Range("H2").Value = 0.05624512345497867
$ws.Range("I2").Value = 0.05624512345497867
$ws.Range("L2").Value = 6.110644904297192
$ws.Range("M2").Value = "[-0.1428067177620278, 12.36409652635641]"
$ws.Range("N2").Value = 0.05523550725474125
$ws.Range("O2").Value = 0.05523550725474125
$ws.Range("P2").Value = -1.320789704211925
$ws.Range("Q2").Value = "[-2.8302636518826962, 0.1886842434588467]"
$ws.Range("R2").Value = 0.08480323047554172
$ws.Range("S2").Value = 0.08480323047554172
$ws.Range("T2").Value = 9.351868281680549
$ws.Range("U2").Value = "[5.683296526364677, 13.02044003699642]"
$ws.Range("V2").Value = [double]"5.874208266032355e-06"
$ws.Range("W2").Value = [double]"5.874208266032355e-06"
$ws.Range("X2").Value = 5.457057057057185
$ws.Range("Y2").Value = -0.7795795795796012
$ws.Range("Z2").Value = 11.69369369369397

# Row 3
$ws.Range("H3").Value = 0.1152084980655497
$ws.Range("I3").Value = 0.1152084980655497
$ws.Range("L3").Value = 5.872170455999512
$ws.Range("M3").Value = "[-1.1328019466853299, 12.877142858684355]"
$ws.Range("N3").Value = 0.09825493480512226
$ws.Range("O3").Value = 0.09825493480512226
$ws.Range("P3").Value = -1.761052938949233
$ws.Range("Q3").Value = "[-3.3963163822592355, -0.12578949563923114]"
$ws.Range("R3").Value = 0.03540268785503486
$ws.Range("S3").Value = 0.03540268785503486
$ws.Range("T3").Value = 11.5260802344408
$ws.Range("U3").Value = "[7.491165945597304, 15.560994523284306]"
$ws.Range("V3").Value = [double]"7.24229488824335e-07"
$ws.Range("W3").Value = [double]"7.24229488824335e-07"
$ws.Range("X3").Value = 7.276076076076251
$ws.Range("Y3").Value = 0.5197197197197347
$ws.Range("Z3").Value = 14.03243243243277

# Row 4
$ws.Range("H4").Value = 0.1491289684539777
$ws.Range("I4").Value = 0.1491289684539777
$ws.Range("L4").Value = 5.643759383735324
$ws.Range("M4").Value = "[-1.8433524704347803, 13.130871237905428]"
$ws.Range("N4").Value = 0.1359527819256268
$ws.Range("O4").Value = 0.1359527819256268
$ws.Range("P4").Value = -2.503210963220696
$ws.Range("Q4").Value = "[-5.6102115055097, 0.6037895790683088]"
$ws.Range("R4").Value = 0.1116398729154915
$ws.Range("S4").Value = 0.1116398729154915
$ws.Range("T4").Value = 10.59015997360164
$ws.Range("U4").Value = "[6.6791495578946645, 14.501170389308623]"
$ws.Range("V4").Value = [double]"2.003848045717405e-06"
$ws.Range("W4").Value = [double]"2.003848045717405e-06"
$ws.Range("X4").Value = 10.34242242242267
$ws.Range("Y4").Value = -2.494654654654715
$ws.Range("Z4").Value = 23.17949949950005

# Row 5
$ws.Range("H5").Value = 0.1201356163724553
$ws.Range("I5").Value = 0.1201356163724553
$ws.Range("L5").Value = 6.527888284871821
$ws.Range("M5").Value = "[-0.5777618838220899, 13.633538453565732]"
$ws.Range("N5").Value = 0.07083707530980954
$ws.Range("O5").Value = 0.07083707530980954
$ws.Range("P5").Value = 3.050395269251351
$ws.Range("Q5").Value = "[1.4402897250691948, 4.660500813433507]"
$ws.Range("R5").Value = 0.0004116590112663854
$ws.Range("S5").Value = 0.0004116590112663854
$ws.Range("T5").Value = 10.42654302362264
$ws.Range("U5").Value = "[6.139549396941618, 14.713536650303652]"
$ws.Range("V5").Value = [double]"1.288507379859638e-05"
$ws.Range("W5").Value = [double]"1.288507379859638e-05"
$ws.Range("X5").Value = 13.35679679679711
$ws.Range("Y5").Value = 6.704384384384543
$ws.Range("Z5").Value = 20.00920920920968

# Row 6
$ws.Range("H6").Value = 0.1334174382042195
$ws.Range("I6").Value = 0.1334174382042195
$ws.Range("L6").Value = 6.957630835678573
$ws.Range("M6").Value = "[-1.132717244475753, 15.0479789158329]"
$ws.Range("N6").Value = 0.0901039642271293
$ws.Range("O6").Value = 0.0901039642271293
$ws.Range("P6").Value = -3.107000542289005
$ws.Range("Q6").Value = "[-4.67936923777939, -1.5346318467986184]"
$ws.Range("R6").Value = 0.0002483342663790467
$ws.Range("S6").Value = 0.0002483342663790467
$ws.Range("T6").Value = 11.45579540038615
$ws.Range("U6").Value = "[6.799210445002728, 16.11238035576957]"
$ws.Range("V6").Value = [double]"1.068620546607058e-05"
$ws.Range("W6").Value = [double]"1.068620546607058e-05"
$ws.Range("X6").Value = 12.83707707707738
$ws.Range("Y6").Value = 6.340580580580737
$ws.Range("Z6").Value = 19.33357357357403

# Row 7
$ws.Range("H7").Value = 0.09864242581245952
$ws.Range("I7").Value = 0.09864242581245952
$ws.Range("L7").Value = 6.443462531136263
$ws.Range("M7").Value = "[-1.555493823902852, 14.442418886175378]"
$ws.Range("N7").Value = 0.1116961060513666
$ws.Range("O7").Value = 0.1116961060513666
$ws.Range("P7").Value = 2.522079387566581
$ws.Range("Q7").Value = "[-0.25786846606042246, 5.302027241193585]"
$ws.Range("R7").Value = 0.0742945526531602
$ws.Range("S7").Value = 0.0742945526531602
$ws.Range("T7").Value = 11.12029586933875
$ws.Range("U7").Value = "[6.904595019179322, 15.33599671949818]"
$ws.Range("V7").Value = [double]"3.224643790211701e-06"
$ws.Range("W7").Value = [double]"3.224643790211701e-06"
$ws.Range("X7").Value = 15.53961961961999
$ws.Range("Y7").Value = 4.053813813813909
$ws.Range("Z7").Value = 27.02542542542606

# Row 8
$ws.Range("H8").Value = 0.2257450165414516
$ws.Range("I8").Value = 0.2257450165414516
$ws.Range("L8").Value = 4.904553319611129
$ws.Range("M8").Value = "[-2.623357426075261, 12.432464065297518]"
$ws.Range("N8").Value = 0.1961001229270845
$ws.Range("O8").Value = 0.1961001229270845
$ws.Range("P8").Value = 1.956026657190042
$ws.Range("Q8").Value = "[-1.1761317842268095, 5.088185098606893]"
$ws.Range("R8").Value = 0.2149508666529132
$ws.Range("S8").Value = 0.2149508666529132
$ws.Range("T8").Value = 10.83843333757162
$ws.Range("U8").Value = "[6.678711567358745, 14.998155107784493]"
$ws.Range("V8").Value = [double]"4.012979385681348e-06"
$ws.Range("W8").Value = [double]"4.012979385681348e-06"
$ws.Range("X8").Value = 17.87835835835878
$ws.Range("Y8").Value = 4.937337337337452
$ws.Range("Z8").Value = 30.81937937938011

# Row 9
$ws.Range("B9").Value = 0
$ws.Range("H9").Value = 0.1633449993693552
$ws.Range("I9").Value = 0.1633449993693552
$ws.Range("L9").Value = 4.955556876336705
$ws.Range("M9").Value = "[-2.1006064472468964, 12.011720199920306]"
$ws.Range("N9").Value = 0.1640949335302511
$ws.Range("O9").Value = 0.1640949335302511
$ws.Range("P9").Value = 2.270500396288119
$ws.Range("Q9").Value = "[-0.8490790955648091, 5.390079888141047]"
$ws.Range("R9").Value = 0.1496261136836932
$ws.Range("S9").Value = 0.1496261136836932
$ws.Range("T9").Value = 8.445075249231856
$ws.Range("U9").Value = "[4.695096647063346, 12.195053851400365]"
$ws.Range("V9").Value = [double]"4.237867527012718e-05"
$ws.Range("W9").Value = [double]"4.237867527012718e-05"
$ws.Range("X9").Value = 16.57905905905945
$ws.Range("Y9").Value = 3.690010010010095
$ws.Range("Z9").Value = 29.46810810810881

# Row 10
$ws.Range("H10").Value = 0.1981819172936363
$ws.Range("I10").Value = 0.1981819172936363
$ws.Range("L10").Value = 5.21181529681497
$ws.Range("M10").Value = "[-2.0641758507567065, 12.487806444386647]"
$ws.Range("N10").Value = 0.1560275104833782
$ws.Range("O10").Value = 0.1560275104833782
$ws.Range("P10").Value = 1.364816027685656
$ws.Range("Q10").Value = "[-1.6918687163476562, 4.421500771718968]"
$ws.Range("R10").Value = 0.3732794149452399
$ws.Range("S10").Value = 0.3732794149452399
$ws.Range("T10").Value = 10.35028317621828
$ws.Range("U10").Value = "[6.195015828994114, 14.505550523442453]"
$ws.Range("V10").Value = [double]"8.695694812876908e-06"
$ws.Range("W10").Value = [double]"8.695694812876908e-06"
$ws.Range("X10").Value = 20.32104104104152
$ws.Range("Y10").Value = 7.691851851852036
$ws.Range("Z10").Value = 32.95023023023101

# Row 11
$ws.Range("H11").Value = 0.04705905970973179
$ws.Range("I11").Value = 0.04705905970973179
$ws.Range("L11").Value = 7.583784839895801
$ws.Range("M11").Value = "[-0.0632042246313631, 15.230773904422964]"
$ws.Range("N11").Value = 0.05184305447463555
$ws.Range("O11").Value = 0.05184305447463555
$ws.Range("P11").Value = 0.4088158608275005
$ws.Range("Q11").Value = "[-1.0817896624973864, 1.8994213841523875]"
$ws.Range("R11").Value = 0.5834138053055788
$ws.Range("S11").Value = 0.5834138053055788
$ws.Range("T11").Value = 11.29777076144854
$ws.Range("U11").Value = "[7.248057968727471, 15.347483554169619]"
$ws.Range("V11").Value = [double]"1.144711346379168e-06"
$ws.Range("W11").Value = [double]"1.144711346379168e-06"
$ws.Range("X11").Value = 24.27091091091149
$ws.Range("Y11").Value = 18.11223223223266
$ws.Range("Z11").Value = 30.42958958959031

# Row 12
$ws.Range("F12").Value = 22.6700000000001
$ws.Range("H12").Value = 0.2759368599990073
$ws.Range("I12").Value = 0.2759368599990073
$ws.Range("L12").Value = 4.613537577063477
$ws.Range("M12").Value = "[-2.8337103135392017, 12.060785467666156]"
$ws.Range("N12").Value = 0.2185846005898524
$ws.Range("O12").Value = 0.2185846005898524
$ws.Range("P12").Value = 0.7484474990534249
$ws.Range("Q12").Value = "[-2.390000417145388, 3.886895415252238]"
$ws.Range("R12").Value = 0.6333286146673089
$ws.Range("S12").Value = 0.6333286146673089
$ws.Range("T12").Value = 10.85473520434803
$ws.Range("U12").Value = "[6.832447141935235, 14.87702326676082]"
$ws.Range("V12").Value = [double]"2.132514482333647e-06"
$ws.Range("W12").Value = [double]"2.132514482333647e-06"
$ws.Range("X12").Value = 19.96956956956966
$ws.Range("Y12").Value = 8.645915915915953
$ws.Range("Z12").Value = 31.29322322322337

# Row 13
$ws.Range("F13").Value = 22.6700000000001
$ws.Range("H13").Value = 0.04014291638903478
$ws.Range("I13").Value = 0.04014291638903478
$ws.Range("L13").Value = 6.630600168642282
$ws.Range("M13").Value = "[0.04786681484450206, 13.213333522440061]"
$ws.Range("N13").Value = 0.04842573817715534
$ws.Range("O13").Value = 0.04842573817715534
$ws.Range("P13").Value = 0.3333421634439633
$ws.Range("Q13").Value = "[-1.1824212590087706, 1.8491055858966972]"
$ws.Range("R13").Value = 0.659934361801958
$ws.Range("S13").Value = 0.659934361801958
$ws.Range("T13").Value = 10.96206625050889
$ws.Range("U13").Value = "[7.311833362918543, 14.612299138099244]"
$ws.Range("V13").Value = [double]"2.644544376817493e-07"
$ws.Range("W13").Value = [double]"2.644544376817493e-07"
$ws.Range("X13").Value = 21.46728728728738
$ws.Range("Y13").Value = 15.99834834834841
$ws.Range("Z13").Value = 26.93622622622635

# Row 14
$ws.Range("F14").Value = 22.6700000000001
$ws.Range("H14").Value = 0.06068021951580738
$ws.Range("I14").Value = 0.06068021951580738
$ws.Range("L14").Value = 7.931837071852719
$ws.Range("M14").Value = "[-0.7980091107540233, 16.66168325445946]"
$ws.Range("N14").Value = 0.07388005733003444
$ws.Range("O14").Value = 0.07388005733003444
$ws.Range("P14").Value = 0.9622896416401172
$ws.Range("Q14").Value = "[-1.440289725069194, 3.3648690083494284]"
$ws.Range("R14").Value = 0.4240848457878572
$ws.Range("S14").Value = 0.4240848457878572
$ws.Range("T14").Value = 13.16908475692914
$ws.Range("U14").Value = "[8.60642779963953, 17.73174171421876]"
$ws.Range("V14").Value = [double]"5.90708456371658e-07"
$ws.Range("W14").Value = [double]"5.90708456371658e-07"
$ws.Range("X14").Value = 19.1980180180181
$ws.Range("Y14").Value = 10.52940940940946
$ws.Range("Z14").Value = 27.86662662662675

# Row 15
$ws.Range("F15").Value = 22.6700000000001
$ws.Range("H15").Value = 0.203687297313258
$ws.Range("I15").Value = 0.203687297313258
$ws.Range("L15").Value = 5.747858128365227
$ws.Range("M15").Value = "[-2.237932373170217, 13.73364862990067]"
$ws.Range("N15").Value = 0.1540839576172495
$ws.Range("O15").Value = 0.1540839576172495
$ws.Range("P15").Value = 1.402552876377426
$ws.Range("Q15").Value = "[-1.6855792415656943, 4.490684994320546]"
$ws.Range("R15").Value = 0.3651956927331435
$ws.Range("S15").Value = 0.3651956927331435
$ws.Range("T15").Value = 10.23465252021963
$ws.Range("U15").Value = "[5.84408119356802, 14.62522384687124]"
$ws.Range("V15").Value = [double]"2.52113971692669e-05"
$ws.Range("W15").Value = [double]"2.52113971692669e-05"
$ws.Range("X15").Value = 17.60952952952961
$ws.Range("Y15").Value = 6.467417417417447
$ws.Range("Z15").Value = 28.75164164164177
